# Updates cryptos price/volume columns (D, E) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Formula = "'26.400.07"
$ws.Cells.Item(2, 5).Formula = "'  +1.15%  "
$ws.Cells.Item(3, 4).Formula = "'1.671.01"
$ws.Cells.Item(3, 5).Formula = "'  +1.09%  "
$ws.Cells.Item(4, 5).Formula = "'  +0.86%  "
$ws.Cells.Item(5, 4).Formula = "'220.68"
$ws.Cells.Item(5, 5).Formula = "'  +1.48%  "
$ws.Cells.Item(6, 4).Formula = "'0.5359"
$ws.Cells.Item(6, 5).Formula = "'  +1.52%  "
$ws.Cells.Item(7, 4).Formula = "'1.010"
$ws.Cells.Item(7, 5).Formula = "'  +0.76%  "
$ws.Cells.Item(8, 4).Formula = "'0.2662"
$ws.Cells.Item(8, 5).Formula = "'  +2.15%  "
$ws.Cells.Item(9, 4).Formula = "'0.06400"
$ws.Cells.Item(9, 5).Formula = "'  +1.36%  "
$ws.Cells.Item(10, 4).Formula = "'21.01"
$ws.Cells.Item(10, 5).Formula = "'  +3.31%  "
$ws.Cells.Item(11, 4).Formula = "'0.07865"
$ws.Cells.Item(11, 5).Formula = "'  +0.95%  "
$ws.Cells.Item(12, 4).Formula = "'4.565"
$ws.Cells.Item(12, 5).Formula = "'  +0.99%  "
$ws.Cells.Item(13, 4).Formula = "'1.664.95"
$ws.Cells.Item(13, 5).Formula = "'  +0.93%  "
$ws.Cells.Item(14, 4).Formula = "'1.899.74"
$ws.Cells.Item(14, 5).Formula = "'  +1.04%  "
$ws.Cells.Item(15, 4).Formula = "'0.5565"
$ws.Cells.Item(15, 5).Formula = "'  +1.60%  "
$ws.Cells.Item(16, 4).Formula = "'0.0" + [char]8325 + "8191"
$ws.Cells.Item(16, 5).Formula = "'  -0.09%  "
$ws.Cells.Item(17, 4).Formula = "'66.24"
$ws.Cells.Item(17, 5).Formula = "'  +1.35%  "
$ws.Cells.Item(18, 4).Formula = "'26.418.07"
$ws.Cells.Item(18, 5).Formula = "'  +1.24%  "
$ws.Cells.Item(19, 4).Formula = "'1.010"
$ws.Cells.Item(19, 5).Formula = "'  +0.79%  "
$ws.Cells.Item(20, 4).Formula = "'4.686"
$ws.Cells.Item(20, 5).Formula = "'  +2.20%  "
$ws.Cells.Item(21, 4).Formula = "'196.25"
$ws.Cells.Item(21, 5).Formula = "'  +2.97%  "
$ws.Cells.Item(22, 4).Formula = "'10.30"
$ws.Cells.Item(22, 5).Formula = "'  +2.40%  "
$ws.Cells.Item(23, 4).Formula = "'6.052"
$ws.Cells.Item(23, 5).Formula = "'  +0.64%  "
$ws.Cells.Item(24, 4).Formula = "'1.011"
$ws.Cells.Item(24, 5).Formula = "'  +0.78%  "
$ws.Cells.Item(25, 4).Formula = "'145.79"
$ws.Cells.Item(25, 5).Formula = "'  +0.74%  "
$ws.Cells.Item(26, 4).Formula = "'0.1228"
$ws.Cells.Item(26, 5).Formula = "'  +0.10%  "
$ws.Cells.Item(27, 4).Formula = "'7.255"
$ws.Cells.Item(28, 4).Formula = "'16.24"
$ws.Cells.Item(28, 5).Formula = "'  +1.72%  "
$ws.Cells.Item(29, 4).Formula = "'1.501"
$ws.Cells.Item(29, 5).Formula = "'  +3.55%  "
$ws.Cells.Item(30, 4).Formula = "'0.05892"
$ws.Cells.Item(30, 5).Formula = "'  +1.99%  "
$ws.Cells.Item(31, 4).Formula = "'1.290"
$ws.Cells.Item(31, 5).Formula = "'  +1.39%  "
$ws.Cells.Item(32, 4).Formula = "'3.583"
$ws.Cells.Item(32, 5).Formula = "'  +0.99%  "
$ws.Cells.Item(33, 4).Formula = "'3.305"
$ws.Cells.Item(33, 5).Formula = "'  +1.32%  "
$ws.Cells.Item(34, 4).Formula = "'1.619"
$ws.Cells.Item(34, 5).Formula = "'  +1.32%  "
$ws.Cells.Item(35, 4).Formula = "'0.9716"
$ws.Cells.Item(35, 5).Formula = "'  +2.71%  "
$ws.Cells.Item(36, 4).Formula = "'2.840"
$ws.Cells.Item(36, 5).Formula = "'  +1.41%  "
$ws.Cells.Item(37, 4).Formula = "'2.431"
$ws.Cells.Item(37, 5).Formula = "'  +0.75%  "
$ws.Cells.Item(38, 4).Formula = "'0.5828"
$ws.Cells.Item(38, 5).Formula = "'  +1.58%  "
$ws.Cells.Item(39, 4).Formula = "'0.01610"
$ws.Cells.Item(39, 5).Formula = "'  -0.10%  "
$ws.Cells.Item(40, 4).Formula = "'1.076.33"
$ws.Cells.Item(40, 5).Formula = "'  +4.33%  "
$ws.Cells.Item(41, 4).Formula = "'0.8666"
$ws.Cells.Item(41, 5).Formula = "'  +1.67%  "
$ws.Cells.Item(42, 4).Formula = "'5.878"
$ws.Cells.Item(42, 5).Formula = "'  +2.90%  "
$ws.Cells.Item(43, 5).Formula = "'  +0.84%  "
$ws.Cells.Item(44, 4).Formula = "'104.32"
$ws.Cells.Item(44, 5).Formula = "'  +0.22%  "
$ws.Cells.Item(45, 4).Formula = "'1.808.02"
$ws.Cells.Item(45, 5).Formula = "'  +0.78%  "
$ws.Cells.Item(46, 4).Formula = "'58.16"
$ws.Cells.Item(47, 4).Formula = "'0.0" + [char]8328 + "106"
$ws.Cells.Item(47, 5).Formula = "'  -5.23%  "
$ws.Cells.Item(48, 4).Formula = "'1.017"
$ws.Cells.Item(48, 5).Formula = "'  +1.64%  "
$ws.Cells.Item(49, 4).Formula = "'0.4398"
$ws.Cells.Item(49, 5).Formula = "'  +1.54%  "
$ws.Cells.Item(50, 4).Formula = "'8.068"
$ws.Cells.Item(50, 5).Formula = "'  +2.57%  "
$ws.Cells.Item(51, 5).Formula = "'  +0.48%  "
